$wb = $excel.ActiveWorkbook

# ============================================================
# Summary sheet
# ============================================================
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Karim Al Ketbi"
$wsSummary.Range("B4").Value = 6385.24
$wsSummary.Range("B6").Value = 704468
$wsSummary.Range("B7").Value = 423711
$wsSummary.Range("B8").Value = 280757
$wsSummary.Range("B9").Value = 1.66

# ============================================================
# Assets sheet
# Insert two new rows (Vehicles: Luxury Car / Mid-range Car) above
# the existing "Liquid Assets" row, and update the totals.
# ============================================================
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Rows("2:3").Insert()

# Copy formatting from the (now shifted) Liquid Assets row down into
# the two freshly inserted rows so borders/fills/number formats match.
$wsAssets.Range("A4:C4").Copy()
$wsAssets.Range("A2:C3").PasteSpecial(-4122)

$wsAssets.Range("A2").Value = "Vehicles"
$wsAssets.Range("B2").Value = "Luxury Car"
$wsAssets.Range("C2").Value = 596199

$wsAssets.Range("A3").Value = "Vehicles"
$wsAssets.Range("B3").Value = "Mid-range Car"
$wsAssets.Range("C3").Value = 98572

$wsAssets.Range("A4").Value = "Liquid Assets"
$wsAssets.Range("B4").Value = "Savings Account"
$wsAssets.Range("C4").Value = 9697

$wsAssets.Range("C5").Value = 704468

# ============================================================
# Liabilities sheet
# Insert two new rows (Auto Loans: Vehicle Loan 1 / Vehicle Loan 2)
# above the existing "Credit Cards" row, and update the totals.
# ============================================================
$wsLiab = $wb.Worksheets.Item("Liabilities")
$wsLiab.Rows("2:3").Insert()

# Copy formatting from the (now shifted) Credit Cards row down into
# the two freshly inserted rows so borders/fills/number formats match.
$wsLiab.Range("A4:E4").Copy()
$wsLiab.Range("A2:E3").PasteSpecial(-4122)

$wsLiab.Range("A2").Value = "Auto Loans"
$wsLiab.Range("B2").Value = "Vehicle Loan 1"
$wsLiab.Range("C2").Value = 357719
$wsLiab.Range("D2").Value = 4968
$wsLiab.Range("E2").Value = 6

$wsLiab.Range("A3").Value = "Auto Loans"
$wsLiab.Range("B3").Value = "Vehicle Loan 2"
$wsLiab.Range("C3").Value = 59143
$wsLiab.Range("D3").Value = 986
$wsLiab.Range("E3").Value = 5

$wsLiab.Range("A4").Value = "Credit Cards"
$wsLiab.Range("B4").Value = "Credit Card Balance"
$wsLiab.Range("C4").Value = 6849
$wsLiab.Range("D4").Value = 342
$wsLiab.Range("E4").Value = 1

$wsLiab.Range("C5").Value = 423711
